# Adds a new column S ("44006") to the accumulated-deaths tracking sheet,
# mirroring the existing columns (header row sum formula, date header,
# per-day counts) and appends a new trailing row (93) for date 44004.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> value for the new column S (rows 3 through 92,
# matching the existing data rows).
$sValues = @{
    3 = 1; 4 = 1; 5 = 1; 6 = 1; 7 = 1; 8 = 1; 9 = 3; 10 = 4;
    11 = 1; 12 = 4; 13 = 6; 14 = 8; 15 = 1; 16 = 5; 17 = 4; 18 = 10; 19 = 9; 20 = 6;
    21 = 7; 22 = 3; 23 = 9; 24 = 4; 25 = 8; 26 = 6; 27 = 5; 28 = 6; 29 = 5; 30 = 10;
    31 = 8; 32 = 5; 33 = 7; 34 = 6; 35 = 7; 36 = 9; 37 = 11; 38 = 7; 39 = 6; 40 = 8;
    41 = 10; 42 = 7; 43 = 19; 44 = 12; 45 = 11; 46 = 14; 47 = 30; 48 = 16; 49 = 22; 50 = 30;
    51 = 19; 52 = 26; 53 = 38; 54 = 35; 55 = 32; 56 = 49; 57 = 46; 58 = 63; 59 = 61; 60 = 65;
    61 = 62; 62 = 70; 63 = 86; 64 = 95; 65 = 86; 66 = 106; 67 = 118; 68 = 113; 69 = 122; 70 = 126;
    71 = 113; 72 = 142; 73 = 133; 74 = 149; 75 = 144; 76 = 148; 77 = 155; 78 = 144; 79 = 135; 80 = 168;
    81 = 154; 82 = 159; 83 = 156; 84 = 164; 85 = 150; 86 = 143; 87 = 159; 88 = 143; 89 = 82; 90 = 40;
    91 = 62; 92 = 78
}

# Write all of the new values first (including the new trailing row 93)
# so the sheet's used range already covers everything before any
# formulas/formatting are applied below.
$ws.Range("S2").Value = 44006
for ($row = 3; $row -le 92; $row++) {
    $ws.Cells.Item($row, 19).Value = $sValues[$row]
}
$ws.Range("A93").Value = 44004
$ws.Range("S93").Value = 17

# Copy formatting (fill/border/number format) from the existing column R
# onto the new column S so the new cells look consistent with the rest of
# the table, for the header row, the date row, and the data rows; also
# copy the row-93 date-column formatting from the row above.
$ws.Range("R1").Copy()
$ws.Range("S1").PasteSpecial(-4122)

$ws.Range("R2").Copy()
$ws.Range("S2").PasteSpecial(-4122)

$ws.Range("R3:R92").Copy()
$ws.Range("S3:S92").PasteSpecial(-4122)

$ws.Range("A92").Copy()
$ws.Range("A93").PasteSpecial(-4122)

$ws.Range("R92").Copy()
$ws.Range("S93").PasteSpecial(-4122)

# Running-total formula for the new column, same pattern as column R.
$ws.Range("S1").Formula = "=SUM(S3:S999)"

# Make sure all formulas (including the new S1 total) are fully
# recalculated with the final values in place.
$excel.Calculate()

# Restore the active selection/scroll position as in the edited workbook.
$ws.Range("U13").Select()
